# Apply updated cryptocurrency Price (D) and Volume(1h) (E) values
# as scraped/updated by the GitHub Actions symbol-list update job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.56"
$ws.Range("E2").Value = "'6.60%"
$ws.Range("D3").Value = "'40.29"
$ws.Range("E3").Value = "'7.25%"
$ws.Range("D4").Value = "'5.277"
$ws.Range("E4").Value = "'2.01%"
$ws.Range("D5").Value = "'0.08099"
$ws.Range("E5").Value = "'2.66%"
$ws.Range("D6").Value = "'8.648"
$ws.Range("E6").Value = "'4.62%"
$ws.Range("D7").Value = "'1.924"
$ws.Range("E7").Value = "'0.78%"
$ws.Range("E8").Value = "'-1.44%"
$ws.Range("D9").Value = "'0.9376"
$ws.Range("E9").Value = "'0.08%"
$ws.Range("D10").Value = "'0.1336"
$ws.Range("E10").Value = "'19.27%"
$ws.Range("D11").Value = "'0.1965"
$ws.Range("E11").Value = "'0.94%"
$ws.Range("D12").Value = "'0.09035"
$ws.Range("E12").Value = "'-0.16%"
$ws.Range("D13").Value = "'0.03472"
$ws.Range("E13").Value = "'4.68%"
$ws.Range("D14").Value = "'0.09572"
$ws.Range("E14").Value = "'-0.31%"
$ws.Range("D15").Value = "'0.001372"
$ws.Range("E15").Value = "'-1.23%"
$ws.Range("D16").Value = "'0.006289"
$ws.Range("E16").Value = "'6.72%"
$ws.Range("D17").Value = "'3.360"
$ws.Range("E17").Value = "'-7.01%"
$ws.Range("D18").Value = "'4.530"
$ws.Range("E18").Value = "'2.43%"
$ws.Range("D19").Value = "'0.3520"
$ws.Range("E19").Value = "'3.19%"
$ws.Range("D20").Value = "'6.506"
$ws.Range("E20").Value = "'0.78%"
$ws.Range("E21").Value = "'3.29%"
$ws.Range("D22").Value = "'0.2573"
$ws.Range("E22").Value = "'2.24%"
$ws.Range("D23").Value = "'0.04444"
$ws.Range("E23").Value = "'0.78%"
$ws.Range("D24").Value = "'0.001227"
$ws.Range("E24").Value = "'-0.55%"
$ws.Range("D25").Value = "'0.004291"
$ws.Range("E25").Value = "'-6.35%"
$ws.Range("D26").Value = "'0.0001296"
$ws.Range("E26").Value = "'-4.84%"
$ws.Range("D27").Value = "'0.0003997"
$ws.Range("E27").Value = "'0.13%"
$ws.Range("D39").Value = "'0.02504"
$ws.Range("E39").Value = "'12.70%"
$ws.Range("D40").Value = "'0.05193"
$ws.Range("E40").Value = "'1.94%"
$ws.Range("D41").Value = "'0.007713"
$ws.Range("E41").Value = "'3.70%"
$ws.Range("D42").Value = "'0.1427"
$ws.Range("E42").Value = "'5.50%"
$ws.Range("D43").Value = "'0.009195"
$ws.Range("E43").Value = "'4.74%"
$ws.Range("D44").Value = "'0.002122"
$ws.Range("E44").Value = "'3.29%"
$ws.Range("D45").Value = "'0.008245"
$ws.Range("E45").Value = "'-4.45%"
$ws.Range("D46").Value = "'0.00006622"
$ws.Range("E46").Value = "'1.35%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.30%"
$ws.Range("D48").Value = "'0.003354"
$ws.Range("E48").Value = "'17.26%"
$ws.Range("E49").Value = "'148.03%"
$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.30%"
$ws.Range("D51").Value = "'0.0002008"
$ws.Range("E51").Value = "'0.30%"
